$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Replace the embedded add-in metadata blob stored in cell A1's comment
$newCommentText = 'Cx0AAB+LCAAAAAAAAAOlWVtvI0kV/istP4GE3d12ZjaJKr3yLVkLO45sh0z2BZW7K3GTcrfpqk7it0UCLVoWIYRm0XIVT4uQGEawKy0zXP7LapIZnvgLnLr0zXZ20mEUTbrOperUqXP5qoLevZ5T45JEzA+DvYpdsyoGCdzQ84PzvUrMz6r248q7Dupeu4Qe4QjPCQdhA7QCtnvN/L3KjPPFrmleXV3Vrhq1MDo365Zlm08G/bE7I3Nc9QPGceCSSqrlvV2r4qC2Nx8Qjj3MsdLcq/TGvVqb+G4HaAMc4HMS1Vox8wPCWDfgPvcJE5oRwZy0O4PvqI059drjmo3MNXom2Yp96im5gqSiazlYlkz8OXHqlr1dtbarDXti7eza1m5jp7a9/ej9RDEVRH3M+JhEl74rCWOO5wupbm03bGvHtho7yNwoBHNlDnDQkHojcukz4rUJpayUR0x9gE2Xw67LOdNCZk5XT/RwEw4ivJhNfE5JOTNGg5YxD7Qt2SQO2g8j4oL/HmTSIbkaRtqtk0UfuJOZH/FlBy9Lz3XMSDRcCCeVU3VQJwx4k5KIHy/grIkHoQAMh0cxQeYdzEyp4zMXvv0gJp5zhinLKxWY6CSMLtgCu+QQ8tgUc1wFNMQeBBz3GffdbNE1BjqKwgXMCIu3Qurtw6xaeAMjnbkXgIvFsq0wvMis28RE8lTl+cKZzjFPxNfoaDwLr4YBXY7jKXMjf0q8TiuR3shDIiG1djtmPJyDFRkJKVqOsoR/kICrZNQhrj/H9IiCE5nTgFkKBNSMeXjm83ZI43nAEptWqOgEdjQh1+kO0zEawuEGwulh0AsSeeXmjayiwii8StdcZ0gn5MhN5qYnss5Zle4ATQfHOkOeiNjlvk+hP+TPIkctRsV4RgjfGBKKg0Qp3Bcdx2lSisxsiCAsIbTBDMeCplGVPxPL2pU/sGzKRt3Au1suYaLDeD6cQvpeyg05NvBWSAi2QFsUBxdAPfH57LCZmL6Bg9SG75Rf5yFI1AXFS0lOnZKnoV7g0tgjKv97wZmMSGGbOpU72WiN1IeUdhAOlpPlAsow83c5fOxVoDHvMh5B6684bhgHPFqKQoFMLfo2HRZPA7kApvfWOYvI92NAHMv9OHDboXf/1TzlnePA5/e3MIwjVf3uryK9JwphzDpElBRZ4++t75bZE4tKic8DMg8D372/t8HJwnrvARthSVbdW4Oo/Lq3PIU2rrqcyPV7q0UAF6GvlVqmyVjo+jJYdXp4OX3zjpTpkDMcU4BqHDrqeVpqV8moyS5WZfIkdBzRpOA5AggzQMKuN6+5gBUE2qu54VwQTACgJ2Nk5uUF4HFJNzjv4+A8BkiR1pVVelpuRTucRDhgYjspglipvJuFUFKnFLJx8iUJmStMNCHzRRhhOgBv+Ps61jQmArgxwHymR9C/KHETz5qZaqpVNCex9m1ishEp20WW69q4QpRCIjIV0M5kMhoSdXkAuUjbmPrTSJXSpF1v4sEpZQAwKbpicyXBYOJ4uGxBh/02WQoEng00XcapnTBU1Irq6YxHW9v1R1ajDuhFjJHc8YhganQhgzkxesElYXwOarvGiDDfgy8f013jPTIlPnQ+6SJ90KW183poPynu0pQm2FukFAUAU5z70DvWBVNOpuCcEhzRZU5QbbUfuiB3++N/3/zm+asXn91+9PTNFz/87z9+9eqfP7t59iP4uP3r324+/qXaphJGEzylRBo0aW1vW40tiLOUhIRzTQl/vdjlknZ6KlFvOkb68iYH7W6vfdBvySKSEhN11UdMcS9chnE2HKtNyIXkkZpJJCgRZ5IUJT0ucHN9yRH3tEtSlM7z71JUvnj98rPXL/98p7Z2WAaw7J2dR1W7/lb8BTdee00uxV/9QuEXwltV61G1Xs8Jr8igkar6qZ96niMu1Fa9YdlpAffSQN4ktMrSM03wubmip0hthYnSEMiPE6YM/AmkSMpWqZAb6BD94idv/vK0IKW9qynFWcA4iV3EYmYykFMfjibGeHg8aneNSXcs4iTj5eTU5F8jrFdP86kQVEEQY/otAzo6dDCjAreeihGeGQS7M2MJmZjLw0KwbaKqhR445aqVB1EYL9SJ5BQy6gbJtJps1NhQayRP+nOt6GSsDeLK1pu/f75JQW+kk6HY9I0jT0MFjiLl+DprP/3Xqy8/fPXixe3zn998+YPCDHqd9OoPcQ7ZlB+mYQ8lT/ebFQo6GUtnXljfzfUXTRS3p6PQDzhz7Mfy4qRHCFRtMZv8jXpzaHlyYukvoK9Q0HuYda+5TmznEJlFAti5wNBtw+x+mRJUDc/8+p/f/u7215/ffvL8zYd/uvnojzcff/L65e/fPPuDyrrbp89vf/pMV/nVRiBtEbdWhfwM+QbiGiIbDdG7ja8++IURhNwAyGHEsiJ99cGnucmEoRKcZDMDjksNKZqwJppXFnpGzpTUhoJeqqIAQFu0sEYqoZtYuPDdbJH3q2IqkXeS8Y3epBozYoSApr4JOykKZ8r31dMqqqUevWPV7brmKmvEFqaY5Vx/QMMpgIyEIR8ZVkQKWl+vkMnK9Q76w1azn4koI4aRRyIRhuoDJZBStJQeS0ZJqOUowAXg58ZUvAutia2z0plzZczU7ytnTU+Uv83vEwUJ1I6jSAGiQL/Xj+MFgOHkGe5uvnyazOHfQ4VV84g4G/c6RT6Mc1xohEW2IEi+LE2apcpUj4m3HAVnD4VrsiHwCs+Z4A79JK+Q1iXgysgUdacbRWG0sfhknERsAEgaKoqZeTyVkWeqULeXnVVCSApe+qGue3qHYYdQwsu9V5uZ9iC8fLAunH1Z1R4bUk87s9zVI3VLNkH+0V4Eyv/7Zq+CrRlFAKzEI1/pR/bktjqCS25Ja9RWpKK4AcLq+hV8348YfyIqgf5SlNOUcqoQ6hNx4VIfcnzqbNUVAQTM/Oxmwcwkdbn600lI+/7cL3kttJL8Lk4CvlwsFITrlYsU0VoOyTUAzNwMUBSn34O2oR5PysymAhZqaaovHiyZfz7jZQ17Z4qJR6ZW1Z2SenXLs7arO4Q0qrYN/2O3XresR+K5U08OlcMnVyUXMZMDy/6k6fwPXaulwwsdAAA='
$ws.Range("A1").Comment.Text($newCommentText)

# 3. Update the custom number format (numFmtId 166) from "0.000" to "###0.000"
#    This format is applied to B27:B43
$ws.Range("B27:B43").NumberFormat = "###0.000"

# 4. Update the label in A11 from "Function Description" to "Function Information"
$ws.Range("A11").Value = "Function Information"
